# Update gh-pages output (北京-漫展信息.xlsx) to match data refreshed at 456a3b4.
# Refreshes "想去人数" (interest count) / "最低票价" (min ticket price) figures
# and marks the "北京· YiYou 运动番only 2.0" event as cancelled across the
# "展览" (row 6), and "全部类型" (row 7) sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1233
$ws.Range("C6").Value = "北京· YiYou 运动番only 2.0（取消）"
$ws.Range("G6").Value = "不可售"
$ws.Range("F7").Value = 7465
$ws.Range("F9").Value = 103
$ws.Range("F10").Value = 2072
$ws.Range("F11").Value = 8119
$ws.Range("F12").Value = 50
$ws.Range("F14").Value = 5558
$ws.Range("F16").Value = 2516
$ws.Range("F17").Value = 1072
$ws.Range("F19").Value = 315
$ws.Range("F21").Value = 85
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = 19.9
$ws.Range("F23").Value = 442
$ws.Range("F24").Value = 1385
$ws.Range("F25").Value = 24
$ws.Range("F26").Value = 2589
$ws.Range("F28").Value = 300
$ws.Range("F29").Value = 105
$ws.Range("F30").Value = 225
$ws.Range("F31").Value = 617
$ws.Range("F34").Value = 1579
$ws.Range("F35").Value = 37
$ws.Range("F37").Value = 2502
$ws.Range("F38").Value = 2252

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 91
$ws.Range("F3").Value = 90
$ws.Range("F5").Value = 31

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1233
$ws.Range("C7").Value = "北京· YiYou 运动番only 2.0（取消）"
$ws.Range("G7").Value = "不可售"
$ws.Range("F8").Value = 7465
$ws.Range("F10").Value = 103
$ws.Range("F11").Value = 2072
$ws.Range("F12").Value = 8119
$ws.Range("F13").Value = 50
$ws.Range("F15").Value = 5558
$ws.Range("F17").Value = 2516
$ws.Range("F18").Value = 1072
$ws.Range("F20").Value = 315
$ws.Range("F22").Value = 85
$ws.Range("F23").Value = 91
$ws.Range("F24").Value = 23
$ws.Range("G24").Value = 19.9
$ws.Range("F25").Value = 90
$ws.Range("F26").Value = 442
$ws.Range("F27").Value = 1385
$ws.Range("F28").Value = 24
$ws.Range("F29").Value = 2589
$ws.Range("F31").Value = 300
$ws.Range("F32").Value = 105
$ws.Range("F33").Value = 225
$ws.Range("F35").Value = 617
$ws.Range("F38").Value = 31
$ws.Range("F39").Value = 1579
$ws.Range("F40").Value = 37
$ws.Range("F42").Value = 2502
$ws.Range("F44").Value = 2252
